$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Bowman Airport station row (row 19): station id, latitude, longitude
$ws.Range("B19").Value = "KBWW"
$ws.Range("E19").Value = 46.1655
$ws.Range("F19").Value = -103.3

# Move the active selection to F19 (as recorded in the saved view state)
$ws.Range("F19").Select() | Out-Null
